# Update cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.388.97"
$ws.Range("D3").Value = "1.856.09"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.84%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.08"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4612"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3713"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07317"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8819"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.96"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07795"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "1.891.79"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.382"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.547"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.78"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009027"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.77"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "27.399.16"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.128"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.52"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "2.116.39"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.936"
$ws.Range("E25").Value = "  +4.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.08"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.106"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.09"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08829"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7783"
$ws.Range("E32").Value = "  +6.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.040"
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.177"
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.504"
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.651"
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01960"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.078"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05231"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.971"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.013"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5152"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1631"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.417"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4821"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.32"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.98"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.650"
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06222"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.78"
$ws.Range("E51").Value = "  +2.03%  "
